# edit.ps1 - applies the diff: capitalize "ranking" -> "Ranking" in two
# specific spots, and append a new "The League Class" section (heading +
# body paragraph) at the end of the document body, before the final
# section break.

$d = $word.ActiveDocument

# --- 1) Heading: The "ranking" Class -> The "Ranking" Class ---------------
$findRange1 = $d.Content
$null = $findRange1.Find.Execute(
    "The “ranking” Class", $false, $false, $false, $false, $false,
    $true, 1, $false, "The “Ranking” Class", 2)

# --- 2) Body: "purpose of the ranking class" -> "purpose of the Ranking class"
$findRange2 = $d.Content
$null = $findRange2.Find.Execute(
    "purpose of the ranking class", $false, $false, $false, $false, $false,
    $true, 1, $false, "purpose of the Ranking class", 2)

# --- 3) Append the new "The League Class" section --------------------------
# Locate the last body paragraph (the one ending "...is outputted. ") and
# insert the new heading + body paragraphs right after it, before sectPr.
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newSectionXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>The “League” Class</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The League class is a based class where different classes can be derived from it. In the league class, there are two inputs that are the league name and number of members of the league. Within the class, there are several virtual functions which each have a different purpose. Each virtual function is called within the base league class with the actual logic for the functions written in the derived class(es).  The first function is the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getLeagueInfo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() virtual function which can be used to get information about the league such as position limits and round limits. A derived class will utilize this function to create the logic for intended operation. The second  function is the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>addMembers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() virtual function which is intended to be used to add league members to a vector container. The next function is the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getLeagueName</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() function which returns the name of the league. Lastly, there are a series of virtual get functions called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getQbLimit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getRbLimit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getWrLimit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>getTeLimit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>andgetRoundLimit</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">() function which are intended to return the values of the quarterback limit, running back limit, wide receiver limit, tight end limit, and round limit respectively. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertPoint.InsertXML($newSectionXml)
